# This workbook is a product-crawl export. Two listings were removed from
# the source crawl (id 6973029 "Severin Tischgrill PG 8565" and id 6995204
# "Electrolux AirFryer Range Explore 6"), and every remaining row's
# timestamp column (O) was refreshed to the new crawl time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTimestamp = "2022-09-15 07:13:45"
$newTimestamp = "2022-09-15 21:02:27"
$idsToRemove = @("6973029", "6995204")

$lastRow = $ws.UsedRange.Rows.Count

# Walk bottom-up so deleting a row doesn't shift the index of rows we still
# need to inspect.
for ($r = $lastRow; $r -ge 2; $r--) {
    $idVal = $ws.Cells.Item($r, 1).Text
    if ($idsToRemove -contains $idVal) {
        $ws.Rows.Item($r).Delete()
    }
}

$lastRow = $ws.UsedRange.Rows.Count
for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 15)
    if ($cell.Text -eq $oldTimestamp) {
        $cell.Value = $newTimestamp
    }
}
